# Apply updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" ('26.091.07')
Set-TextValue "E2" ('  -0.24%  ')

Set-TextValue "D3" ('1.646.15')
Set-TextValue "E3" ('  -0.68%  ')

Set-TextValue "D4" ('1.013')
Set-TextValue "E4" ('  +0.38%  ')

Set-TextValue "D5" ('215.99')
Set-TextValue "E5" ('  -0.98%  ')

Set-TextValue "D6" ('0.5046')
Set-TextValue "E6" ('  -2.04%  ')

Set-TextValue "D7" ('1.013')
Set-TextValue "E7" ('  +0.44%  ')

Set-TextValue "D8" ('0.2585')
Set-TextValue "E8" ('  +0.69%  ')

Set-TextValue "D9" ('0.06442')
Set-TextValue "E9" ('  +0.32%  ')

Set-TextValue "E10" ('  -1.71%  ')

Set-TextValue "D11" ('0.07747')
Set-TextValue "E11" ('  -0.56%  ')

Set-TextValue "D12" ('1.646.98')
Set-TextValue "E12" ('  -0.77%  ')

Set-TextValue "D13" ('4.258')
Set-TextValue "E13" ('  -1.06%  ')

Set-TextValue "D14" ('1.874.81')
Set-TextValue "E14" ('  -0.56%  ')

Set-TextValue "D15" ('0.5465')
Set-TextValue "E15" ('  -1.50%  ')

Set-TextValue "D16" ('0.0' + [string][char]8325 + '7939')
Set-TextValue "E16" ('  -1.30%  ')

Set-TextValue "D17" ('63.77')
Set-TextValue "E17" ('  -0.96%  ')

Set-TextValue "D18" ('26.104.84')
Set-TextValue "E18" ('  -0.34%  ')

Set-TextValue "E19" ('  +0.54%  ')

Set-TextValue "D20" ('203.86')
Set-TextValue "E20" ('  -3.20%  ')

Set-TextValue "D21" ('4.312')
Set-TextValue "E21" ('  -1.99%  ')

Set-TextValue "D22" ('10.02')
Set-TextValue "E22" ('  -0.53%  ')

Set-TextValue "D23" ('5.976')
Set-TextValue "E23" ('  +1.54%  ')

Set-TextValue "D24" ('1.014')
Set-TextValue "E24" ('  +0.54%  ')

Set-TextValue "D25" ('1.956')
Set-TextValue "E25" ('  +11.28%  ')

Set-TextValue "D26" ('142.21')
Set-TextValue "E26" ('  -1.25%  ')

Set-TextValue "D27" ('0.1158')
Set-TextValue "E27" ('  -0.20%  ')

Set-TextValue "D28" ('15.70')
Set-TextValue "E28" ('  -0.39%  ')

Set-TextValue "D29" ('6.753')
Set-TextValue "E29" ('  -3.00%  ')

Set-TextValue "D30" ('0.05070')
Set-TextValue "E30" ('  -3.63%  ')

Set-TextValue "D31" ('1.246')
Set-TextValue "E31" ('  -0.74%  ')

Set-TextValue "D32" ('3.266')
Set-TextValue "E32" ('  -3.11%  ')

Set-TextValue "D33" ('3.203')
Set-TextValue "E33" ('  -0.35%  ')

Set-TextValue "D34" ('1.548')
Set-TextValue "E34" ('  -1.30%  ')

Set-TextValue "D35" ('2.352')
Set-TextValue "E35" ('  -0.72%  ')

Set-TextValue "D36" ('0.8995')
Set-TextValue "E36" ('  -2.63%  ')

Set-TextValue "D37" ('2.624')
Set-TextValue "E37" ('  -4.62%  ')

Set-TextValue "D38" ('0.5651')
Set-TextValue "E38" ('  -1.34%  ')

Set-TextValue "D39" ('1.154.12')
Set-TextValue "E39" ('  -0.90%  ')

Set-TextValue "D40" ('0.01576')
Set-TextValue "E40" ('  -0.95%  ')

Set-TextValue "D41" ('2.580')
Set-TextValue "E41" ('  +0.20%  ')

Set-TextValue "D43" ('5.672')
Set-TextValue "E43" ('  +0.12%  ')

Set-TextValue "D44" ('0.8173')
Set-TextValue "E44" ('  -2.90%  ')

Set-TextValue "D45" ('100.02')
Set-TextValue "E45" ('  +0.08%  ')

Set-TextValue "D46" ('1.786.18')
Set-TextValue "E46" ('  -0.47%  ')

Set-TextValue "D47" ('0.0' + [string][char]8328 + '115')
Set-TextValue "E47" ('  +5.36%  ')

Set-TextValue "D48" ('0.4548')
Set-TextValue "E48" ('  +1.00%  ')

Set-TextValue "D49" ('1.013')
Set-TextValue "E49" ('  +0.30%  ')

Set-TextValue "D50" ('55.08')
Set-TextValue "E50" ('  -1.59%  ')

Set-TextValue "D51" ('0.05047')
Set-TextValue "E51" ('  -0.91%  ')
